# Added logging and output file: new "Unsolvable" sheet with an unsolvable
# Sudoku puzzle, conditional formatting highlighting the empty (0) cells,
# and make it the active tab.

$wb = $excel.ActiveWorkbook

# --- create the new worksheet -------------------------------------------
$ws = $wb.Worksheets.Add()
$ws.Name = "Unsolvable"

# --- fill in the (unsolvable) sudoku grid --------------------------------
$data = @(
  @(0,5,0,0,2,0,0,1,0),
  @(0,6,0,1,0,4,8,9,0),
  @(0,1,0,0,0,0,2,7,0),
  @(1,0,6,9,0,0,0,3,0),
  @(5,0,0,0,0,0,9,0,1),
  @(9,8,0,0,1,3,0,4,0),
  @(0,7,1,2,0,5,4,6,0),
  @(0,0,5,6,0,1,0,0,0),
  @(6,0,0,8,0,0,1,0,0)
)

for ($r = 0; $r -lt $data.Count; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
  }
}

# --- conditional formatting: highlight zeros (unsolved cells) -----------
$rng = $ws.Range("A1:I9")
$cf = $rng.FormatConditions.Add(1, 3, "0")
$cf.Font.Color = 0x0006009C
$cf.Interior.Color = 0x00CEC7FF

# move it to the end of the tab strip (after the last existing sheet)
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- make the new sheet the active tab, with the same selection that
#     was present in the authored workbook --------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Select()
$last.Range("M10").Select() | Out-Null
